# Fruta / hortaliza, semanal
# Re-assign D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) values across rows 2-24
# (row 6 stays untouched) following the shuffled weekly snapshot order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> source row (the row whose original D/J/K/L/M/P values should end up
# in this row after the edit).
$mapping = @{
    2  = 7
    3  = 14
    4  = 23
    5  = 13
    7  = 9
    8  = 4
    9  = 15
    10 = 17
    11 = 10
    12 = 19
    13 = 21
    14 = 3
    15 = 24
    16 = 5
    17 = 22
    18 = 16
    19 = 2
    20 = 12
    21 = 18
    22 = 11
    23 = 8
    24 = 20
}

# Snapshot the original values for the affected columns before writing
# anything back, since this is a permutation of existing row data. Use
# Value2 for reading since Value can surface a Variant wrapper instead of
# the underlying primitive in this runtime.
$original = @{}
foreach ($row in $mapping.Keys) {
    $original[$row] = @{
        D = $ws.Cells.Item($row, 4).Value2
        J = $ws.Cells.Item($row, 10).Value2
        K = $ws.Cells.Item($row, 11).Value2
        L = $ws.Cells.Item($row, 12).Value2
        M = $ws.Cells.Item($row, 13).Value2
        P = $ws.Cells.Item($row, 16).Value2
    }
}

foreach ($row in $mapping.Keys) {
    $src = $original[$mapping[$row]]
    $ws.Cells.Item($row, 4).Value = $src.D
    $ws.Cells.Item($row, 10).Value = $src.J
    $ws.Cells.Item($row, 11).Value = $src.K
    $ws.Cells.Item($row, 12).Value = $src.L
    $ws.Cells.Item($row, 13).Value = $src.M
    $ws.Cells.Item($row, 16).Value = $src.P
}
